# Update the dSF column (F) values for several rows, per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -3
$ws.Range("F4").Value  = -7
$ws.Range("F5").Value  = -6
$ws.Range("F6").Value  = -3
$ws.Range("F10").Value = -8
$ws.Range("F13").Value = 14
$ws.Range("F14").Value = 0
$ws.Range("F16").Value = -6
